# Update plots for each sample
#
# The underlying analysis was re-run and produced updated results for the
# CYP2D6_10B / CYP2D6_002 marker (sample S1):
#   - the mutant-allele peak height threshold moved from 1000 to 800
#   - a peak that previously went undetected is now detected, with its
#     measured peak/size/height/status populated and the failure message
#     cleared
#   - the marker's called genotype moved from homozygous wildtype (CC) to
#     heterozygous (CT)
#   - the sample-level overall genotype call updated accordingly

$wb = $excel.ActiveWorkbook

# --- peak_table: mutant peak height for CYP2D6_10B (row 3) ---
$wsPeak = $wb.Worksheets.Item("peak_table")
$wsPeak.Range("O3").Value = 800

# --- allele_table: detection results for CYP2D6_10B / S1 (row 5) ---
$wsAllele = $wb.Worksheets.Item("allele_table")
$wsAllele.Range("K5").Value = 800
$wsAllele.Range("M5").Value = $true
$wsAllele.Range("N5").Value = 32
$wsAllele.Range("O5").Value = 34.32
$wsAllele.Range("P5").Value = 949
$wsAllele.Range("Q5").Value = "ok"
$wsAllele.Range("R5").Value = ""

# --- marker_table: genotype/phenotype call for CYP2D6_10B (row 3) ---
$wsMarker = $wb.Worksheets.Item("marker_table")
$wsMarker.Range("G3").Value = "CT"
$wsMarker.Range("H3").Value = "heterozygous"

# --- genotype_result: overall sample genotype (row 2) ---
$wsGenotype = $wb.Worksheets.Item("genotype_result")
$wsGenotype.Range("B2").Value = "*1/*10B"
